$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 932.2017155010354
$ws.Range("E2").Value = 1143.41330678195
$ws.Range("F2").Value = 703.940969505364
$ws.Range("G2").Value = 703.9409696685137
